$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A12").Value = "HE_meningitis_Campign"
$ws.Range("A13").Value = "HE_meningitis_CampignSize"
$ws.Range("A11").Value = "HE_meningitis_Routine"
$ws.Range("A3").Value = "HE_cholera_Campign"
$ws.Range("A4").Value = "HE_cholera_CampignSize"
$ws.Range("A35").Value = "HE_yellowfever_Campign"
$ws.Range("A34").Value = "HE_yellowfever_CampignSize"
$ws.Range("A36").Value = "HE_yellowfever_Routine"
$ws.Range("A17").Value = "HE_polio_Routine"
$ws.Range("A18").Value = "HE_surviving_infants"
$ws.Range("A9").Value = "HE_measles_Routine"

$ws.Range("B3").Select()
